# Apply "multi rnn and xavier init" results to the resultsTable sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 ("15. 5way multiCell + xavier +100dim"): fill in train acc and test acc/epoch note
$ws.Range("B19").Value = 56.14
$ws.Range("C19").Value = "5/10"

# Row 20 ("16. 5way xavier"): fill in train acc and test acc/epoch note
$ws.Range("B20").Value = 54.46
$ws.Range("C20").Value = "5/20"

# Move the active selection to D20, matching the last cell touched during editing.
$ws.Range("D20").Select()
